# Added RMSE and Manhattan evaluation to instance-based vote count
# prediction (only InceptionV3 and color hist included).
#
# The "BASADO EN INSTANCIAS" (instance-based) table, at rows 12-20, gets two
# new k columns (k=15, k=20) added to its right-hand sub-table (columns
# H:N, headed by the "RESULTADOS DE DISTANCIA MANHATTAN" merged title) and
# the HIST-EUCLID/HIST-COS/INCV3-EUCLID/INCV3-COS/k=9 figures that were
# previously missing for that sub-table are filled in as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers on row 14 (right-hand, "instance based" table) ---
$ws.Range("L14").Value = "k=9"
$ws.Range("M14").Value = "k=15"
$ws.Range("N14").Value = "k=20"

# --- Row 15 (HIST-EUCLID) ---
$ws.Range("H15").Value = 6.03
$ws.Range("I15").Value = 4.93
$ws.Range("J15").Value = 4.66
$ws.Range("K15").Value = 4.54
$ws.Range("L15").Value = 4.51
$ws.Range("M15").Value = 4.44
$ws.Range("N15").Value = 4.43

# --- Row 16 (HIST-COS) ---
$ws.Range("H16").Value = 6.08
$ws.Range("I16").Value = 4.92
$ws.Range("J16").Value = 4.71
$ws.Range("K16").Value = 4.57
$ws.Range("L16").Value = 4.52
$ws.Range("M16").Value = 4.43
$ws.Range("N16").Value = 4.39

# --- Row 19 (INCV3-EUCLID) ---
$ws.Range("H19").Value = 5.36
$ws.Range("I19").Value = 4.32
$ws.Range("J19").Value = 4.13
$ws.Range("K19").Value = 4.07
$ws.Range("L19").Value = 3.98
$ws.Range("M19").Value = 4.05
$ws.Range("N19").Value = 4.07

# --- Row 20 (INCV3-COS) ---
$ws.Range("H20").Value = 5.6
$ws.Range("I20").Value = 4.4
$ws.Range("J20").Value = 4.24
$ws.Range("K20").Value = 4.14
$ws.Range("L20").Value = 4.08
$ws.Range("M20").Value = 4.02

# N20 explicitly carries an applied "General" number format in the source
# workbook (applyNumberFormat="1"), distinguishing it from its neighbours.
$ws.Range("N20").NumberFormat = "General"
$ws.Range("N20").Value = 4.07

# --- View/selection bookkeeping: author scrolled down/right and left the
# selection on N17 before saving ---
$ws.Range("N17").Select()
